# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D) for the fea4cfaa... row (row 5)
# on both the zh-cn and de-de language sheets to reflect a new handoff report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-14 15:08:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-14 15:09:14"
